# emx/dist/rd3_portal_novelomics.xlsx
# Commit: "added secondary ID column"
#
# 1. Bump the package description version/date string for rd3_portal_novelomics.
# 2. Add a new "CNAG_barcode" attribute (bool) for the
#    rd3_portal_novelomics_shipment entity, inserted right before the
#    existing "processed" attribute row.

$wb = $excel.ActiveWorkbook

# --- 1. Update package description (sheet "packages", row 3 = rd3_portal_novelomics) ---
$wsPackages = $wb.Worksheets.Item("packages")
$wsPackages.Range("C3").Value = "Staging tables for novel omics sample and experiment metadata (v1.1.0, 2021-09-28)"

# --- 2. Insert new attribute row in sheet "attributes" ---
$wsAttr = $wb.Worksheets.Item("attributes")

# "processed" currently sits at row 56, directly above "molgenis_id" (row 57),
# both belonging to rd3_portal_novelomics_shipment. Insert a new blank row
# above "processed" to hold the new CNAG_barcode attribute.
$wsAttr.Rows.Item(56).Insert()

$wsAttr.Range("A56").Value = "rd3_portal_novelomics_shipment"
$wsAttr.Range("B56").Value = "CNAG_barcode"
$wsAttr.Range("C56").Value = "FALSE"
$wsAttr.Range("D56").Value = "bool"
$wsAttr.Range("E56").Value = "FALSE"
$wsAttr.Range("F56").Value = "TRUE"
